$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. SampTiss: drop the obsolete tissue-type values (Abnormal, Non-neoplastic,
#    Unavailable, Unspecified), keeping Not Reported / Normal / Peritumoral /
#    Tumor / Unknown. Delete from the bottom up so row numbers of
#    not-yet-deleted rows stay stable.
# ---------------------------------------------------------------------------
$sampTiss = $wb.Worksheets.Item("SampTiss")
$sampTiss.Activate()
$sampTiss.Rows.Item(9).Delete()   # Unspecified
$sampTiss.Rows.Item(7).Delete()   # Unavailable
$sampTiss.Rows.Item(6).Delete()   # Non-neoplastic
$sampTiss.Rows.Item(2).Delete()   # Abnormal
[void]$sampTiss.Range("B14").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "SampTumGrd" sheet (Sample Tumor Grade) right after
#    "SampTumor" with its allowed values, as a new Federation test case.
# ---------------------------------------------------------------------------
$sampTumor = $wb.Worksheets.Item("SampTumor")
$sampTumGrd = $wb.Worksheets.Add($null, $sampTumor)
$sampTumGrd.Name = "SampTumGrd"

$gradeValues = @(
    "G1 Low Grade",
    "G2 Intermediate Grade",
    "G3 High Grade",
    "G4 Anaplastic",
    "GB Borderline",
    "GX Grade Cannot Be Assessed",
    "Not Applicable",
    "Not Reported",
    "Unknown"
)
for ($i = 0; $i -lt $gradeValues.Count; $i++) {
    $sampTumGrd.Cells.Item($i + 1, 1).Value = $gradeValues[$i]
}

$sampTumGrd.Columns.Item(1).AutoFit()
$sampTumGrd.Activate()
[void]$sampTumGrd.Range("D8").Select()
